$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff: (cell reference -> new value)
$updates = @(
    @{ Cell = 'D2'; Value = '27.327.35' }
    @{ Cell = 'E2'; Value = '  +3.19%  ' }
    @{ Cell = 'D3'; Value = '1.836.88' }
    @{ Cell = 'E3'; Value = '  +3.33%  ' }
    @{ Cell = 'E4'; Value = '  +3.32%  ' }
    @{ Cell = 'D5'; Value = '317.04' }
    @{ Cell = 'E5'; Value = '  +3.23%  ' }
    @{ Cell = 'D6'; Value = '1.019' }
    @{ Cell = 'E6'; Value = '  +1.49%  ' }
    @{ Cell = 'D7'; Value = '0.4338' }
    @{ Cell = 'E7'; Value = '  +0.91%  ' }
    @{ Cell = 'D8'; Value = '0.3715' }
    @{ Cell = 'E8'; Value = '  +1.55%  ' }
    @{ Cell = 'D9'; Value = '0.07326' }
    @{ Cell = 'E9'; Value = '  +1.92%  ' }
    @{ Cell = 'E10'; Value = '  +2.85%  ' }
    @{ Cell = 'D11'; Value = '2.097.23' }
    @{ Cell = 'E11'; Value = '  +17.04%  ' }
    @{ Cell = 'E12'; Value = '  +4.31%  ' }
    @{ Cell = 'D14'; Value = '6.675' }
    @{ Cell = 'E14'; Value = '  +3.25%  ' }
    @{ Cell = 'D15'; Value = '0.07118' }
    @{ Cell = 'E15'; Value = '  +3.27%  ' }
    @{ Cell = 'D16'; Value = '82.01' }
    @{ Cell = 'E16'; Value = '  +3.33%  ' }
    @{ Cell = 'D17'; Value = '1.025' }
    @{ Cell = 'E17'; Value = '  +1.67%  ' }
    @{ Cell = 'D18'; Value = '0.000009005' }
    @{ Cell = 'E18'; Value = '  +3.24%  ' }
    @{ Cell = 'D19'; Value = '1.016' }
    @{ Cell = 'E19'; Value = '  +1.25%  ' }
    @{ Cell = 'D20'; Value = '15.42' }
    @{ Cell = 'E20'; Value = '  +2.14%  ' }
    @{ Cell = 'D21'; Value = '27.368.69' }
    @{ Cell = 'E21'; Value = '  +3.33%  ' }
    @{ Cell = 'D22'; Value = '5.232' }
    @{ Cell = 'E22'; Value = '  +2.00%  ' }
    @{ Cell = 'D23'; Value = '11.10' }
    @{ Cell = 'E23'; Value = '  -0.37%  ' }
    @{ Cell = 'D24'; Value = '2.311.95' }
    @{ Cell = 'E24'; Value = '  +15.41%  ' }
    @{ Cell = 'D25'; Value = '156.26' }
    @{ Cell = 'E25'; Value = '  +2.49%  ' }
    @{ Cell = 'D26'; Value = '1.899' }
    @{ Cell = 'E26'; Value = '  +1.60%  ' }
    @{ Cell = 'D27'; Value = '18.54' }
    @{ Cell = 'E27'; Value = '  +2.38%  ' }
    @{ Cell = 'D28'; Value = '5.279' }
    @{ Cell = 'E28'; Value = '  +3.14%  ' }
    @{ Cell = 'E29'; Value = '  +7.10%  ' }
    @{ Cell = 'D30'; Value = '115.37' }
    @{ Cell = 'E30'; Value = '  +0.86%  ' }
    @{ Cell = 'D31'; Value = '0.09004' }
    @{ Cell = 'E31'; Value = '  +0.37%  ' }
    @{ Cell = 'D32'; Value = '1.200' }
    @{ Cell = 'D33'; Value = '0.7591' }
    @{ Cell = 'E33'; Value = '  +3.72%  ' }
    @{ Cell = 'D34'; Value = '4.459' }
    @{ Cell = 'E34'; Value = '  +2.61%  ' }
    @{ Cell = 'D35'; Value = '2.839' }
    @{ Cell = 'E35'; Value = '  +3.70%  ' }
    @{ Cell = 'D36'; Value = '1.020' }
    @{ Cell = 'E36'; Value = '  +1.48%  ' }
    @{ Cell = 'D37'; Value = '1.147' }
    @{ Cell = 'E37'; Value = '  +5.65%  ' }
    @{ Cell = 'B38'; Value = 'VeChain' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D38'; Value = '0.01951' }
    @{ Cell = 'E38'; Value = '  +3.16%  ' }
    @{ Cell = 'B39'; Value = 'Hedera' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D39'; Value = '0.05259' }
    @{ Cell = 'E39'; Value = '  +1.56%  ' }
    @{ Cell = 'D40'; Value = '0.5158' }
    @{ Cell = 'E40'; Value = '  +4.43%  ' }
    @{ Cell = 'D41'; Value = '2.787' }
    @{ Cell = 'E41'; Value = '  +7.52%  ' }
    @{ Cell = 'D42'; Value = '0.1662' }
    @{ Cell = 'E42'; Value = '  +2.78%  ' }
    @{ Cell = 'D43'; Value = '6.525' }
    @{ Cell = 'E43'; Value = '  +3.31%  ' }
    @{ Cell = 'D44'; Value = '8.446' }
    @{ Cell = 'E44'; Value = '  +4.82%  ' }
    @{ Cell = 'D45'; Value = '107.84' }
    @{ Cell = 'E45'; Value = '  +2.63%  ' }
    @{ Cell = 'D46'; Value = '10.52' }
    @{ Cell = 'E46'; Value = '  +3.38%  ' }
    @{ Cell = 'D47'; Value = '1.021' }
    @{ Cell = 'E47'; Value = '  +1.68%  ' }
    @{ Cell = 'B48'; Value = 'Decentraland' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D48'; Value = '0.4629' }
    @{ Cell = 'E48'; Value = '  +2.74%  ' }
    @{ Cell = 'B49'; Value = 'RenderToken' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D49'; Value = '1.904' }
    @{ Cell = 'E49'; Value = '  +9.28%  ' }
    @{ Cell = 'D50'; Value = '1.662' }
    @{ Cell = 'E50'; Value = '  +3.10%  ' }
    @{ Cell = 'D51'; Value = '0.06274' }
    @{ Cell = 'E51'; Value = '  +1.17%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. '1.200', '27.327.35')
    # are not reinterpreted by Excel as numbers/dates.
    $range.NumberFormat = '@'
    $range.Value = $u.Value
}

